$p = $ppt.ActivePresentation
try {
    $v = $p.HasNotesMaster
    Write-Output "HasNotesMaster=[$v]"
} catch {
    Write-Output "ERROR: $_"
}
